$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: Shift 2 -> Shift 3
$ws.Name = "Shift 3 - 2025-02-16"

# Move the "RESPALDO DE BASE DE DATOS" activity (currently the first data
# row) to the end of the activity log, so the other activities move up one
# row. Copy row 2 down below the last row, then delete the original row 2
# (shifting rows 3-5 up into 2-4).
$ws.Range("A2:E2").Copy($ws.Range("A6:E6"))
$ws.Rows.Item(2).Delete()

# Update Shift number for every activity row: 2 -> 3
$ws.Range("B2:B5").Value = 3

# Filter the Engineer column down to only the engineer(s) actually working
# the shift, dropping the extra "administrador" / "Usuario de BC" names.
$engineer = "Renato Hacel Cal y Mayor Rodríguez"
$ws.Range("E2:E5").Value = $engineer
